# This script updates the Icam4-Itga2b NATMI LR-pair output sheet with
# recomputed TPM-based values (per commit message: "update scripts wuth new tpm").
# The underlying ligand/receptor expression numbers were regenerated, which
# changes derived columns (expression cells/rate, average/total expression,
# specificity scores and edge weights) across data rows 2-10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2564746666666666
$ws.Range("H2").Value = 0.7694239999999999
$ws.Range("I2").Value = 0.1818007399394835
$ws.Range("J2").Value = 0.1818007399394835
$ws.Range("M2").Value = 0.3883076666666667
$ws.Range("N2").Value = 1.164923
$ws.Range("O2").Value = 0.1188638477168776
$ws.Range("P2").Value = 0.1188638477168776
$ws.Range("Q2").Value = 0.09959107937244444
$ws.Range("R2").Value = 0.8963197143519999
$ws.Range("S2").Value = 0.02160953546698244
$ws.Range("T2").Value = 0.02160953546698244
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2564746666666666
$ws.Range("H3").Value = 0.7694239999999999
$ws.Range("I3").Value = 0.1818007399394835
$ws.Range("J3").Value = 0.1818007399394835
$ws.Range("O3").Value = 0.6829215134520935
$ws.Range("P3").Value = 0.6829215134520935
$ws.Range("Q3").Value = 0.572191561671111
$ws.Range("R3").Value = 5.149724055039998
$ws.Range("S3").Value = 0.1241556364661825
$ws.Range("T3").Value = 0.1241556364661826
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2564746666666666
$ws.Range("H4").Value = 0.7694239999999999
$ws.Range("I4").Value = 0.1818007399394835
$ws.Range("J4").Value = 0.1818007399394835
$ws.Range("O4").Value = 0.1982146388310289
$ws.Range("P4").Value = 0.1982146388310289
$ws.Range("Q4").Value = 0.1660758103306667
$ws.Range("R4").Value = 1.494682292976
$ws.Range("S4").Value = 0.03603556800631854
$ws.Range("T4").Value = 0.03603556800631855
$ws.Range("I5").Value = 0.7694380609030022
$ws.Range("J5").Value = 0.7694380609030022
$ws.Range("M5").Value = 0.3883076666666667
$ws.Range("N5").Value = 1.164923
$ws.Range("O5").Value = 0.1188638477168776
$ws.Range("P5").Value = 0.1188638477168776
$ws.Range("Q5").Value = 0.4215008531927778
$ws.Range("R5").Value = 3.793507678735
$ws.Range("S5").Value = 0.09145836849874404
$ws.Range("T5").Value = 0.09145836849874404
$ws.Range("I6").Value = 0.7694380609030022
$ws.Range("J6").Value = 0.7694380609030022
$ws.Range("O6").Value = 0.6829215134520935
$ws.Range("P6").Value = 0.6829215134520935
$ws.Range("S6").Value = 0.5254658050595223
$ws.Range("T6").Value = 0.5254658050595223
$ws.Range("I7").Value = 0.7694380609030022
$ws.Range("J7").Value = 0.7694380609030022
$ws.Range("O7").Value = 0.1982146388310289
$ws.Range("P7").Value = 0.1982146388310289
$ws.Range("S7").Value = 0.1525138873447358
$ws.Range("T7").Value = 0.1525138873447358
$ws.Range("G8").Value = 0.06878966666666667
$ws.Range("I8").Value = 0.0487611991575143
$ws.Range("J8").Value = 0.0487611991575143
$ws.Range("M8").Value = 0.3883076666666667
$ws.Range("N8").Value = 1.164923
$ws.Range("O8").Value = 0.1188638477168776
$ws.Range("P8").Value = 0.1188638477168776
$ws.Range("Q8").Value = 0.02671155495411111
$ws.Range("R8").Value = 0.240403994587
$ws.Range("S8").Value = 0.00579594375115112
$ws.Range("T8").Value = 0.00579594375115112
$ws.Range("G9").Value = 0.06878966666666667
$ws.Range("I9").Value = 0.0487611991575143
$ws.Range("J9").Value = 0.0487611991575143
$ws.Range("O9").Value = 0.6829215134520935
$ws.Range("P9").Value = 0.6829215134520935
$ws.Range("Q9").Value = 0.1534688291377778
$ws.Range("S9").Value = 0.03330007192638861
$ws.Range("T9").Value = 0.03330007192638861
$ws.Range("G10").Value = 0.06878966666666667
$ws.Range("I10").Value = 0.0487611991575143
$ws.Range("J10").Value = 0.0487611991575143
$ws.Range("O10").Value = 0.1982146388310289
$ws.Range("P10").Value = 0.1982146388310289
$ws.Range("Q10").Value = 0.04454357922566667
$ws.Range("R10").Value = 0.400892213031
$ws.Range("S10").Value = 0.00966518347997457
$ws.Range("T10").Value = 0.00966518347997457
